$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order row (row 5): SKU, Item, Quantity, Cost Per, Total Cost
# Force text format so numeric-looking values are stored as text (matching
# the existing rows, which are all stored as strings) rather than being
# auto-converted to numbers by Excel.
$row = $ws.Range("A5:E5")
$row.NumberFormat = "@"

$ws.Range("A5").Value = "173339"
$ws.Range("B5").Value = "Chobani - Drinkable Yogurt"
$ws.Range("C5").Value = "6"
$ws.Range("D5").Value = "17.99"
$ws.Range("E5").Value = "107.94"
